# Generate Report for Handoff
# Rename the source file GUID from 5001f62d-6ee4-4523-b8ed-be061910262c
# to 6d7ed314-7f6d-476e-8a1f-399549ae3259 across all sheets, update the
# associated xliff hash tokens, and refresh the handoff timestamps.

$wb = $excel.ActiveWorkbook

$oldGuid = "5001f62d-6ee4-4523-b8ed-be061910262c"
$newGuid = "6d7ed314-7f6d-476e-8a1f-399549ae3259"

$oldZhHash = "8fb95be92108e22f158c86233ffd1aeff5dda161"
$newZhHash = "bdb88c07381ff71f0d9af58254a7fe932a1ec477"

$oldDeHash = "8fb95be92108e22f158c86233ffd1aeff5dda161"
$newDeHash = "bdb88c07381ff71f0d9af58254a7fe932a1ec477"

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Hyperlinks.Item(1).TextToDisplay = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-08-18 22:59:07"

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("A2").Value = "$newGuid.md"
$wsZh.Hyperlinks.Item(1).TextToDisplay = "$newGuid.md"
$wsZh.Range("G2").Value = "$newGuid.$newZhHash.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-08-18 22:58:58"

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("A2").Value = "$newGuid.md"
$wsDe.Hyperlinks.Item(1).TextToDisplay = "$newGuid.md"
$wsDe.Range("G2").Value = "$newGuid.$newDeHash.de-de.xlf"
$wsDe.Range("H2").Value = "2016-08-18 22:59:07"
